# Auto-generated edit script for 杭州-漫展信息.xlsx refresh
$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" (Exhibitions): refresh F (interest count) / G (min price/status) ----
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 1322
$wsExpo.Range("F3").Value = 2176
$wsExpo.Range("G3").Value = "70"
$wsExpo.Range("F4").Value = 439
$wsExpo.Range("F5").Value = 196
$wsExpo.Range("G5").Value = "已停售"
$wsExpo.Range("F6").Value = 439
$wsExpo.Range("G6").Value = "130"
$wsExpo.Range("F7").Value = 65
$wsExpo.Range("G7").Value = "已停售"
$wsExpo.Range("F8").Value = 568
$wsExpo.Range("F9").Value = 115
$wsExpo.Range("F10").Value = 193
$wsExpo.Range("F11").Value = 840
$wsExpo.Range("F12").Value = 70
$wsExpo.Range("F14").Value = 7
$wsExpo.Range("F15").Value = 4816
$wsExpo.Range("F16").Value = 2732
$wsExpo.Range("F17").Value = 870
$wsExpo.Range("F18").Value = 655
$wsExpo.Range("F19").Value = 367
$wsExpo.Range("F21").Value = 747
$wsExpo.Range("F22").Value = 1648
$wsExpo.Range("F23").Value = 64
$wsExpo.Range("F24").Value = 692
$wsExpo.Range("F25").Value = 314
$wsExpo.Range("F26").Value = 116
$wsExpo.Range("F27").Value = 217

# ---- Sheet "演出" (Performances): drop the cancelled SACG row, renumber the index column ----
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Rows.Item(2).Delete()
$wsShow.Range("A2").Value = 1
$wsShow.Range("A3").Value = 2
$wsShow.Range("A4").Value = 3

# ---- Sheet "本地生活" (Local life): refresh F (interest count) ----
$wsLocal = $wb.Worksheets.Item("本地生活")
$wsLocal.Range("F2").Value = 140

# ---- Sheet "全部类型" (All types): drop the same cancelled SACG row, renumber, refresh F/G ----
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Rows.Item(3).Delete()
for ($r = 3; $r -le 31; $r++) {
    $wsAll.Cells.Item($r, 1).Value = $r - 1
}
$wsAll.Range("F2").Value = 140
$wsAll.Range("F4").Value = 1322
$wsAll.Range("F5").Value = 2176
$wsAll.Range("G5").Value = "70"
$wsAll.Range("F6").Value = 439
$wsAll.Range("F7").Value = 196
$wsAll.Range("G7").Value = "已停售"
$wsAll.Range("F8").Value = 439
$wsAll.Range("G8").Value = "130"
$wsAll.Range("F9").Value = 65
$wsAll.Range("G9").Value = "已停售"
$wsAll.Range("F10").Value = 568
$wsAll.Range("F11").Value = 115
$wsAll.Range("F12").Value = 193
$wsAll.Range("F13").Value = 840
$wsAll.Range("F14").Value = 70
$wsAll.Range("F18").Value = 7
$wsAll.Range("F19").Value = 4816
$wsAll.Range("F20").Value = 2732
$wsAll.Range("F21").Value = 870
$wsAll.Range("F22").Value = 655
$wsAll.Range("F23").Value = 367
$wsAll.Range("F25").Value = 747
$wsAll.Range("F26").Value = 1648
$wsAll.Range("F27").Value = 64
$wsAll.Range("F28").Value = 692
$wsAll.Range("F29").Value = 314
$wsAll.Range("F30").Value = 116
$wsAll.Range("F31").Value = 217

Write-Output "done"
